# Append a new row (row 47) of sensor data to each of the four worksheets,
# mirroring the pattern already used by the preceding rows.

$wb = $excel.ActiveWorkbook

# Data for the new row on each worksheet, in worksheet order
# (ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER, ROW02-MID-LIFTER).
$newRows = @(
    @{
        A = "2025-03-06 06:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        A = "2025-03-06 06:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        A = "2025-03-06 06:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        A = "2025-03-06 06:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $rowData = $newRows[$i - 1]

    $targetRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($targetRow, 1).Value = $rowData.A
    $ws.Cells.Item($targetRow, 2).Value = $rowData.B
    $ws.Cells.Item($targetRow, 3).Value = $rowData.C
    $ws.Cells.Item($targetRow, 4).Value = $rowData.D
    $ws.Cells.Item($targetRow, 5).Value = $rowData.E
    $ws.Cells.Item($targetRow, 6).Value = $rowData.F

    # Column G holds a 24-digit identifier that exceeds floating point
    # precision, so it must be written (and stay) as text.
    $gCell = $ws.Cells.Item($targetRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $rowData.G
    $gCell.Style = "Normal"

    $ws.Cells.Item($targetRow, 8).Value = $rowData.H
    $ws.Cells.Item($targetRow, 9).Value = $rowData.I
}
